$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Responsables")
$ws.Columns("AY").Insert()
